# Remove the two (nearly) empty sheets that aren't part of the final workbook.
$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item("Sheet1").Delete() | Out-Null
$wb.Worksheets.Item("Sheet5").Delete() | Out-Null

# Rename the three remaining sheets to their new, descriptive names.
$wsText = $wb.Worksheets.Item("Sheet2")
$wsText.Name = "text_box"

$wsCheck = $wb.Worksheets.Item("Sheet3")
$wsCheck.Name = "check_box"

$wsRadio = $wb.Worksheets.Item("Sheet4")
$wsRadio.Name = "radio_box"

# Scroll "text_box" so row 5 is at the top of the view (selection stays at F3).
$wsText.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1

# Scroll "check_box" so row 2 is at the top of the view (selection stays at A2).
$wsCheck.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1

# "radio_box" becomes the active sheet, scrolled to row 2, with the
# selection moved from B5 to C10.
$wsRadio.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$wsRadio.Range("C10").Select()
